$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.918.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -5.21%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.104.82'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -5.64%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '557.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '162.30'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -9.76%  '

$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.587'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.65%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.097.10'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.70%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.72'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.80%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.116'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -7.80%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.376'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.33%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.635.31'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.04%  '

$ws.Range("E14").Value = '  -1.90%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '63.035.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.00%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '24.55'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -7.58%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.114.17'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000152'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '400.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -6.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.36'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.77%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.61%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.87%  '

$ws.Range("E23").Value = '  +0.39%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.69'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.47%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.69'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.42%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.197'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.49%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.476'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.19%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000101'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -11.62%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.65'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.80%  '

$ws.Range("E30").Value = '  -0.32%  '

$ws.Range("E31").Value = '  -0.05%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.77'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.52%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.85'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.80'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -6.55%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.18'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.65%  '

$ws.Range("B36").Value = 'Fetch.AI'
$ws.Range("C36").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.09'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.24%  '

$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '151.65'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.32'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -7.59%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.697.51'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.85%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.64'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -8.27%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '23.39'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -10.38%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.02'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.70%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '38.12'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.99%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.690'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.76%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0602'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -8.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0253'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.44%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.17'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -12.52%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '282.05'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -9.11%  '

$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '20.72'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.03%  '

$ws.Range("B50").Value = 'FirstDigitalUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.06%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0970'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.35%  '
